$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 69, pushing the existing rows 69:186 down to 70:187.
$ws.Rows("69:69").Insert()

# Populate the newly-inserted row 69 with the new weekly observation.
$ws.Range("A69").Value = 8
$ws.Range("B69").Value = "Terminal La Palmera de La Serena"
$ws.Range("C69").Value = "Coquimbo"
$ws.Range("D69").Value = 45082
$ws.Range("E69").Value = 4
$ws.Range("F69").Value = 100112052
$ws.Range("G69").Value = "Albahaca"
$ws.Range("H69").Value = "Sin especificar"
$ws.Range("I69").Value = "Primera"
$ws.Range("J69").Value = 800
$ws.Range("K69").Value = 2800
$ws.Range("L69").Value = 3000
$ws.Range("M69").Value = 2900
$ws.Range("N69").Value = "$/paquete"
$ws.Range("O69").Value = "Región de Arica y Parinacota"
$ws.Range("P69").Value = 2900
$ws.Range("Q69").Value = 1
$ws.Range("R69").Value = "Hortaliza"
